$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "[Botha's lark] (https://www.birdlife.org.za/red-list/bothas-lark/)"
$ws.Range("A3").Value = "[Blue swallow] (https://www.birdlife.org.za/red-list/blue-swallow/)"
$ws.Range("A4").Value = "[White-backed vulture] (https://www.birdlife.org.za/red-list/white-backed-vulture/)"
$ws.Range("A5").Value = "[Lappet-faced vulture] (https://www.birdlife.org.za/red-list/lappet-faced-vulture/)"
$ws.Range("A6").Value = "[Hooded vulture] (https://www.birdlife.org.za/red-list/lappet-faced-vulture/)"
$ws.Range("A7").Value = "[Bearded vulture] (https://www.birdlife.org.za/red-list/bearded-vulture/)"
$ws.Range("A8").Value = "[African penguin] (https://www.birdlife.org.za/red-list/african-penguin/)"
$ws.Range("A9").Value = "[Bank cormorant] (https://www.birdlife.org.za/red-list/bank-cormorant/)"
$ws.Range("A10").Value = "[Cape gannet] (https://www.birdlife.org.za/red-list/cape-gannet/)"
$ws.Range("A11").Value = "[White-headed vulture] (https://www.birdlife.org.za/red-list/white-headed-vulture/)"
$ws.Range("A12").Value = "[Black harrier] (https://www.birdlife.org.za/red-list/black-harrier/)"
$ws.Range("A13").Value = "[Marion island breeding spp] (https://www.birdlife.org.za/?s=Marion+Island+breeding+spp&id=323007&post_type=red-list)"

$ws.Range("A14").Select() | Out-Null
